# Append row 10 to the results table (MLP run: number_of_seasons, 100 trees).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(10, 1).Value = 10
$ws.Cells.Item(10, 2).Value = "number_of_seasons"
$ws.Cells.Item(10, 3).Value = 0.763870780420232
$ws.Cells.Item(10, 4).Value = 0.2072117123926988
$ws.Cells.Item(10, 5).Value = 100

# F10/G10 stay blank (like the other data rows) but must remain present as
# empty text cells rather than being cleared outright. A leading apostrophe
# forces an explicit (empty) text entry; resetting the style back to Normal
# afterwards drops the transient "quote prefix" formatting it implies, so
# the cell ends up identical to the sheet's other blank text cells.
$ws.Cells.Item(10, 6).Value = "'"
$ws.Cells.Item(10, 6).Style = "Normal"
$ws.Cells.Item(10, 7).Value = "'"
$ws.Cells.Item(10, 7).Style = "Normal"
